$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# The sheet currently has data rows 2..42 (row 1 = header).
# Old rows 40, 41, 42 (ids 149, 150, 151) describe three timeline
# entries dated 44224 / 44227 / 44229.
#
# The new commit adds a fresh COVID patient timeline (reported
# 3 Feb 2021) whose locations ("ตลาดบางแค", "ซีคอน บางแค") need to be
# inserted chronologically *before* the existing 44224 entry, plus one
# more "unspecified_location" stop that lands right after it (before
# the 44227 entry). Net effect: 3 new rows, old rows shift down by
# 2 then by 1 more.
#
# Final layout (by old identity):
#   row 40 -> NEW  (ตลาดบางแค)
#   row 41 -> NEW  (ซีคอน บางแค)
#   row 42 -> was old row 40 (ตลาดนัดสายไหม, 44224)
#   row 43 -> NEW  (โรงงานแถว ถ.เพชรเกษม เขตภาษีเจริญ, unspecified_location)
#   row 44 -> was old row 41 (polli's cafe, 44227)
#   row 45 -> was old row 42 (มหาวิทยาลัยเกษตรศาสตร์ กำแพงแสน, 44229)
# ------------------------------------------------------------------

# Insert the three new rows at the right spots (old row indices shift
# as each insert happens, so we do the ones nearer the top first).
$ws.Rows.Item(40).Insert()
$ws.Rows.Item(40).Insert()
$ws.Rows.Item(43).Insert()

# ------------------------------------------------------------------
# Fill in the brand-new rows.
# ------------------------------------------------------------------

# Row 40 - ตลาดบางแค
$ws.Range("A40").Value = 149
$ws.Range("B40").Value = "patient_has_been_here"
$ws.Range("C40").Value = 44220
$ws.Range("D40").Value = "ตลาดบางแค"
$ws.Range("E40").Value = 13.71189925
$ws.Range("F40").Value = 100.4275498
$ws.Range("G40").Value = 0
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = " กทม.เปิดไทม์ไลน์ ผู้ป่วยโควิด  ประจำวันที่  3 กพ."
$ws.Range("M40").Value = "https://www.facebook.com/earthpongsakornk/posts/456692252407342"
$ws.Range("N40").Value = 0

# Row 41 - ซีคอน บางแค
$ws.Range("A41").Value = 150
$ws.Range("B41").Value = "patient_has_been_here"
$ws.Range("C41").Value = 44220
$ws.Range("D41").Value = "ซีคอน บางแค"
$ws.Range("E41").Value = 13.71171772
$ws.Range("F41").Value = 100.4340445
$ws.Range("G41").Value = 0
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = " กทม.เปิดไทม์ไลน์ ผู้ป่วยโควิด  ประจำวันที่  3 กพ."
$ws.Range("M41").Value = "https://www.facebook.com/earthpongsakornk/posts/456692252407342"
$ws.Range("N41").Value = 0

# Row 42 (was old row 40) - renumber id only, rest of the data carried
# over automatically by the row insert.
$ws.Range("A42").Value = 151

# Row 43 - โรงงานแถว ถ.เพชรเกษม เขตภาษีเจริญ (unspecified_location)
$ws.Range("A43").Value = 152
$ws.Range("B43").Value = "unspecified_location"
$ws.Range("C43").Value = 44224
$ws.Range("D43").Value = "โรงงานแถว ถ.เพชรเกษม เขตภาษีเจริญ"
$ws.Range("E43").Value = 13.71657541
$ws.Range("F43").Value = 100.4430568
$ws.Range("G43").Value = 0
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = " กทม.เปิดไทม์ไลน์ ผู้ป่วยโควิด  ประจำวันที่  3 กพ."
$ws.Range("M43").Value = "https://www.facebook.com/earthpongsakornk/posts/456692252407342"
$ws.Range("N43").Value = 0

# Row 44 (was old row 41) - renumber id only.
$ws.Range("A44").Value = 153

# Row 45 (was old row 42) - renumber id only.
$ws.Range("A45").Value = 154

# ------------------------------------------------------------------
# Hyperlinks: row inserts don't reliably keep the M-column hyperlink
# references in sync, so rebuild the whole collection from the
# current M-column text top to bottom.
# ------------------------------------------------------------------
$ws.Hyperlinks.Delete()

$lastRow = 45
$r = 2
while ($r -le $lastRow) {
    $cell = $ws.Range("M" + $r)
    $ws.Hyperlinks.Add($cell, $cell.Value)
    $r = $r + 1
}
